# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# Column G (header "K", row 1) holds the newly-regenerated "K" values
# that replace the previous "Strike#" derived numbers for each data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> new K value (column G), per the regenerated save_data
$s_vals = @{
    2  = 5
    3  = 5
    4  = 5
    5  = 4
    6  = 6
    7  = 3
    8  = 3
    9  = 3
    10 = 4
    11 = 7
    12 = 2
    13 = 2
    14 = 0
}

foreach ($row in $s_vals.Keys) {
    $ws.Cells.Item($row, 7).Value = $s_vals[$row]
}
